$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-14 06:48:56'
$ws.Range('G2').Value = '122 cm'
$ws.Range('I2').Value = '10.8 mm'
$ws.Range('N2').Value = '-1.2 °C 6:08 TU'
$ws.Range('E3').Value = '2026-02-14 06:48:58'
$ws.Range('I3').Value = '4.8 mm'
$ws.Range('N3').Value = '-5.4 °C 6:27 TU'
$ws.Range('E4').Value = '2026-02-14 06:49:01'
$ws.Range('H4').NumberFormat = '@'
$ws.Range('H4').Value = '84%'
$ws.Range('H4').NumberFormat = 'General'
$ws.Range('J4').Value = '990.0 hPa'
$ws.Range('N4').Value = '6.1 °C 6:22 TU'
$ws.Range('O4').Value = '8.2 °C'
$ws.Range('E5').Value = '2026-02-14 06:49:04'
$ws.Range('G5').Value = '119 cm'
$ws.Range('I5').Value = '7.9 mm'
$ws.Range('N5').Value = '-5.2 °C 6:14 TU'
$ws.Range('O5').Value = '-4.8 °C'
$ws.Range('E6').Value = '2026-02-14 06:49:06'
$ws.Range('J6').Value = '990.0 hPa'
$ws.Range('N6').Value = '6.3 °C 6:08 TU'
$ws.Range('E7').Value = '2026-02-14 06:49:09'
$ws.Range('J7').Value = '990.3 hPa'
$ws.Range('E8').Value = '2026-02-14 06:49:12'
$ws.Range('J8').Value = '989.9 hPa'
$ws.Range('N8').Value = '6.5 °C 6:11 TU'
$ws.Range('E9').Value = '2026-02-14 06:49:15'
$ws.Range('E10').Value = '2026-02-14 06:49:17'
$ws.Range('N10').Value = '6.1 °C 6:26 TU'
$ws.Range('O10').Value = '6.7 °C'
$ws.Range('E11').Value = '2026-02-14 06:49:20'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '89%'
$ws.Range('H11').NumberFormat = 'General'
$ws.Range('M11').Value = '7.6 °C 6:26 TU'
$ws.Range('O11').Value = '3.6 °C'
$ws.Range('E12').Value = '2026-02-14 06:49:23'
$ws.Range('H12').NumberFormat = '@'
$ws.Range('H12').Value = '60%'
$ws.Range('H12').NumberFormat = 'General'
$ws.Range('E13').Value = '2026-02-14 06:49:25'
$ws.Range('H13').NumberFormat = '@'
$ws.Range('H13').Value = '85%'
$ws.Range('H13').NumberFormat = 'General'
$ws.Range('J13').Value = '991.6 hPa'
$ws.Range('O13').Value = '2.3 °C'
$ws.Range('E14').Value = '2026-02-14 06:49:28'
$ws.Range('E15').Value = '2026-02-14 06:49:31'
$ws.Range('H15').NumberFormat = '@'
$ws.Range('H15').Value = '62%'
$ws.Range('H15').NumberFormat = 'General'
$ws.Range('E16').Value = '2026-02-14 06:49:33'
$ws.Range('I16').Value = '5.0 mm'
$ws.Range('O16').Value = '-5.6 °C'
$ws.Range('E17').Value = '2026-02-14 06:49:36'
$ws.Range('N17').Value = '-0.1 °C 6:22 TU'
$ws.Range('E18').Value = '2026-02-14 06:49:39'
$ws.Range('J18').Value = '990.2 hPa'
$ws.Range('N18').Value = '6.3 °C 6:26 TU'
$ws.Range('O18').Value = '7.3 °C'
$ws.Range('E19').Value = '2026-02-14 06:49:41'
$ws.Range('N19').Value = '2.4 °C 6:16 TU'
$ws.Range('E20').Value = '2026-02-14 06:49:44'
$ws.Range('I20').Value = '1.5 mm'
$ws.Range('N20').Value = '-5.5 °C 6:28 TU'
$ws.Range('E21').Value = '2026-02-14 06:49:47'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '94%'
$ws.Range('H21').NumberFormat = 'General'
$ws.Range('J21').Value = '992.6 hPa'
$ws.Range('O21').Value = '1.4 °C'
$ws.Range('E22').Value = '2026-02-14 06:49:50'
$ws.Range('H22').NumberFormat = '@'
$ws.Range('H22').Value = '89%'
$ws.Range('H22').NumberFormat = 'General'
$ws.Range('N22').Value = '-6.9 °C 6:23 TU'
$ws.Range('E23').Value = '2026-02-14 06:49:52'
$ws.Range('I23').Value = '11.9 mm'
$ws.Range('N23').Value = '-5.9 °C 6:29 TU'
$ws.Range('O23').Value = '-5.4 °C'
$ws.Range('E24').Value = '2026-02-14 06:49:55'
$ws.Range('J24').Value = '994.0 hPa'
$ws.Range('E25').Value = '2026-02-14 06:49:58'
$ws.Range('I25').Value = '21.7 mm'
$ws.Range('O25').Value = '-4.0 °C'
$ws.Range('E26').Value = '2026-02-14 06:50:01'
$ws.Range('E27').Value = '2026-02-14 06:50:03'
$ws.Range('N27').Value = '-3.6 °C 6:02 TU'
$ws.Range('E28').Value = '2026-02-14 06:50:06'
$ws.Range('J28').Value = '990.5 hPa'
$ws.Range('N28').Value = '4.3 °C 6:19 TU'
$ws.Range('O28').Value = '5.1 °C'
$ws.Range('E29').Value = '2026-02-14 06:50:08'
$ws.Range('E30').Value = '2026-02-14 06:50:11'
$ws.Range('H30').NumberFormat = '@'
$ws.Range('H30').Value = '71%'
$ws.Range('H30').NumberFormat = 'General'
$ws.Range('J30').Value = '989.5 hPa'
$ws.Range('O30').Value = '10.3 °C'
$ws.Range('E31').Value = '2026-02-14 06:50:14'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '77%'
$ws.Range('H31').NumberFormat = 'General'
$ws.Range('J31').Value = '989.0 hPa'
$ws.Range('N31').Value = '8.5 °C 6:23 TU'
$ws.Range('E32').Value = '2026-02-14 06:50:17'
$ws.Range('E33').Value = '2026-02-14 06:50:19'
$ws.Range('H33').NumberFormat = '@'
$ws.Range('H33').Value = '68%'
$ws.Range('H33').NumberFormat = 'General'
$ws.Range('J33').Value = '990.3 hPa'
$ws.Range('E34').Value = '2026-02-14 06:50:22'
$ws.Range('H34').NumberFormat = '@'
$ws.Range('H34').Value = '74%'
$ws.Range('H34').NumberFormat = 'General'
$ws.Range('I34').Value = '0.2 mm'
$ws.Range('E35').Value = '2026-02-14 06:50:25'
$ws.Range('J35').Value = '995.8 hPa'
$ws.Range('N35').Value = '1.7 °C 6:25 TU'
$ws.Range('O35').Value = '2.5 °C'
$ws.Range('E36').Value = '2026-02-14 06:50:27'
$ws.Range('H36').NumberFormat = '@'
$ws.Range('H36').Value = '63%'
$ws.Range('H36').NumberFormat = 'General'
$ws.Range('J36').Value = '989.8 hPa'
$ws.Range('E37').Value = '2026-02-14 06:50:30'
$ws.Range('J37').Value = '991.6 hPa'
$ws.Range('N37').Value = '3.0 °C 6:03 TU'
$ws.Range('E38').Value = '2026-02-14 06:50:33'
$ws.Range('L38').Value = '10.8 km/h - 297º 6:16 TU'
$ws.Range('N38').Value = '6.4 °C 6:00 TU'
$ws.Range('E39').Value = '2026-02-14 06:50:36'
$ws.Range('I39').Value = '3.7 mm'
$ws.Range('N39').Value = '-5.9 °C 6:28 TU'
$ws.Range('E40').Value = '2026-02-14 06:50:38'
$ws.Range('H40').NumberFormat = '@'
$ws.Range('H40').Value = '96%'
$ws.Range('H40').NumberFormat = 'General'
$ws.Range('J40').Value = '993.1 hPa'
$ws.Range('M40').Value = '8.2 °C 6:12 TU'
$ws.Range('O40').Value = '2.8 °C'
$ws.Range('E41').Value = '2026-02-14 06:50:41'
$ws.Range('J41').Value = '991.5 hPa'
$ws.Range('N41').Value = '10.7 °C 6:29 TU'
$ws.Range('O41').Value = '11.5 °C'
$ws.Range('E42').Value = '2026-02-14 06:50:44'
$ws.Range('H42').NumberFormat = '@'
$ws.Range('H42').Value = '82%'
$ws.Range('H42').NumberFormat = 'General'
$ws.Range('N42').Value = '8.6 °C 6:03 TU'
$ws.Range('E43').Value = '2026-02-14 06:50:46'
$ws.Range('H43').NumberFormat = '@'
$ws.Range('H43').Value = '93%'
$ws.Range('H43').NumberFormat = 'General'
$ws.Range('N43').Value = '4.1 °C 6:27 TU'
$ws.Range('O43').Value = '5.7 °C'
$ws.Range('E44').Value = '2026-02-14 06:50:49'
$ws.Range('I44').Value = '15.2 mm'
$ws.Range('N44').Value = '-5.5 °C 6:21 TU'
$ws.Range('E45').Value = '2026-02-14 06:50:52'
$ws.Range('I45').Value = '4.9 mm'
$ws.Range('J45').Value = '997.2 hPa'
$ws.Range('E46').Value = '2026-02-14 06:50:54'
$ws.Range('J46').Value = '995.1 hPa'
$ws.Range('N46').Value = '9.2 °C 6:28 TU'
$ws.Range('O46').Value = '10.2 °C'
